# Update column F ("dSF") values for specific rows as per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    4  = -6
    6  = 1
    10 = -2
    11 = -1
    14 = 0
    19 = -2
    20 = 2
    26 = -1
    31 = 0
    35 = 1
    40 = 2
    45 = -1
    46 = -2
    47 = -1
    51 = 0
    52 = -1
    53 = -2
    59 = 1
    60 = 1
    62 = -2
    71 = -5
    72 = 3
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
